# Saldo_guide.xlsx update
# - Refresh "Dt. Referencia" (column G) for every data row from 26/04/2024 (45408)
#   to 29/04/2024 (45411).
# - A handful of rows had their "Saldo Previsto" (D), "Vl. Projetado" (E) and
#   "Vl. Total" (H) values recalculated / corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference date for every data row (2 through 310) in one shot.
$ws.Range("G2:G310").Value = 45411

# Row-specific corrections to Saldo Previsto / Vl. Projetado / Vl. Total.
$ws.Range("D12").Value = 973.73
$ws.Range("E12").Value = 0
$ws.Range("H12").Value = 973.73

$ws.Range("D44").Value = 60.67
$ws.Range("H44").Value = 60.67

$ws.Range("D124").Value = 551.16
$ws.Range("H124").Value = 551.16

$ws.Range("D129").Value = 108.69
$ws.Range("H129").Value = 108.69

$ws.Range("D184").Value = 608.55999999999995
$ws.Range("E184").Value = 0
$ws.Range("H184").Value = 608.55999999999995

$ws.Range("D217").Value = 0.04
$ws.Range("H217").Value = 0.04
